# Apply benchmark update per commit: 2026-02-20 07:05:36 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "33,33 TL - 33,33 TL"

# Row 3
$ws.Range("G3").Value = ""

# Row 4
$ws.Range("G4").Value = ""

# Row 5
$ws.Range("G5").Value = ""

# Row 6
$ws.Range("C6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("J6").Value = ""

# Row 7
$ws.Range("D7").Value = "%1,6"

# Row 8
$ws.Range("G8").Value = ""

# Row 9
$ws.Range("G9").Value = ""

# Row 10
$ws.Range("G10").Value = ""

# Row 11
$ws.Range("G11").Value = ""

# Row 12
$ws.Range("C12").Value = ""
$ws.Range("G12").Value = ""

# Row 13
$ws.Range("C13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 69,62 TL"

# Row 14
$ws.Range("C14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("J14").Value = ""
